# Apply updated team-specific time-matrix probabilities to Sheet1
# (Hampton_B.xlsx team_specific_matrix data). Values are recomputed
# fractions (e.g. x/27 style ratios) replacing earlier rounded figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1814814814814815
$ws.Range("C2").Value = 0.5740740740740741
$ws.Range("J2").Value = 0.01111111111111111
$ws.Range("P2").Value = 0.1666666666666667
$ws.Range("S2").Value = 0.06666666666666667
$ws.Range("B3").Value = 0.02380952380952381
$ws.Range("C3").Value = 0.06547619047619048
$ws.Range("J3").Value = 0.02976190476190476
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.130952380952381
$ws.Range("J4").Value = 0.02857142857142857
$ws.Range("P4").Value = 0.7714285714285715
$ws.Range("S4").Value = 0.2
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.05333333333333334
$ws.Range("D6").Value = 0.01333333333333333
$ws.Range("F6").Value = 0.03555555555555556
$ws.Range("J6").Value = 0.32
$ws.Range("O6").Value = 0.02222222222222222
$ws.Range("Q6").Value = 0.1377777777777778
$ws.Range("R6").Value = 0.08444444444444445
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.08888888888888889
$ws.Range("D7").Value = 0.02222222222222222
$ws.Range("F7").Value = 0.06111111111111111
$ws.Range("J7").Value = 0.1333333333333333
$ws.Range("Q7").Value = 0.15
$ws.Range("R7").Value = 0.09444444444444444
$ws.Range("S7").Value = 0.45
$ws.Range("B8").Value = 0.08798283261802575
$ws.Range("D8").Value = 0.01716738197424893
$ws.Range("E8").Value = 0.002145922746781116
$ws.Range("F8").Value = 0.05579399141630902
$ws.Range("J8").Value = 0.1244635193133047
$ws.Range("O8").Value = 0.02575107296137339
$ws.Range("Q8").Value = 0.1802575107296137
$ws.Range("R8").Value = 0.07939914163090128
$ws.Range("S8").Value = 0.4270386266094421
$ws.Range("B9").Value = 0.08205128205128205
$ws.Range("D9").Value = 0.01538461538461539
$ws.Range("E9").Value = 0.005128205128205128
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1487179487179487
$ws.Range("O9").Value = 0.01538461538461539
$ws.Range("Q9").Value = 0.1897435897435897
$ws.Range("R9").Value = 0.09230769230769231
$ws.Range("S9").Value = 0.3846153846153846
$ws.Range("B10").Value = 0.0935672514619883
$ws.Range("D10").Value = 0.01388888888888889
$ws.Range("E10").Value = 0.002923976608187134
$ws.Range("F10").Value = 0.06140350877192982
$ws.Range("J10").Value = 0.1198830409356725
$ws.Range("O10").Value = 0.01973684210526316
$ws.Range("Q10").Value = 0.2412280701754386
$ws.Range("R10").Value = 0.09137426900584796
$ws.Range("S10").Value = 0.3559941520467836
$ws.Range("G11").Value = 0.1654929577464789
$ws.Range("J11").Value = 0.0880281690140845
$ws.Range("K11").Value = 0.2112676056338028
$ws.Range("L11").Value = 0.5352112676056338
$ws.Range("G12").Value = 0.7320261437908496
$ws.Range("J12").Value = 0.2222222222222222
$ws.Range("K12").Value = 0.0130718954248366
$ws.Range("L12").Value = 0.006535947712418301
$ws.Range("S12").Value = 0.0261437908496732
$ws.Range("F13").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.6046511627906976
$ws.Range("J13").Value = 0.3255813953488372
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.004761904761904762
$ws.Range("H15").Value = 0.1809523809523809
$ws.Range("I15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.3571428571428572
$ws.Range("K15").Value = 0.0380952380952381
$ws.Range("O15").Value = 0.02857142857142857
$ws.Range("S15").Value = 0.3238095238095238
$ws.Range("F16").Value = 0.03125
$ws.Range("H16").Value = 0.1354166666666667
$ws.Range("I16").Value = 0.07291666666666667
$ws.Range("J16").Value = 0.46875
$ws.Range("K16").Value = 0.078125
$ws.Range("M16").Value = 0.02083333333333333
$ws.Range("O16").Value = 0.06770833333333333
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.01785714285714286
$ws.Range("H17").Value = 0.1547619047619048
$ws.Range("I17").Value = 0.1031746031746032
$ws.Range("J17").Value = 0.4325396825396826
$ws.Range("K17").Value = 0.08531746031746032
$ws.Range("M17").Value = 0.00992063492063492
$ws.Range("O17").Value = 0.05357142857142857
$ws.Range("S17").Value = 0.1428571428571428
$ws.Range("F18").Value = 0.05069124423963134
$ws.Range("H18").Value = 0.1474654377880184
$ws.Range("I18").Value = 0.09216589861751152
$ws.Range("J18").Value = 0.4147465437788018
$ws.Range("K18").Value = 0.1059907834101382
$ws.Range("M18").Value = 0.01382488479262673
$ws.Range("N18").Value = 0.004608294930875576
$ws.Range("O18").Value = 0.05990783410138249
$ws.Range("S18").Value = 0.1105990783410138
$ws.Range("F19").Value = 0.02313030069390902
$ws.Range("H19").Value = 0.2274479568234387
$ws.Range("I19").Value = 0.07401696222050887
$ws.Range("J19").Value = 0.3677717810331534
$ws.Range("K19").Value = 0.10254433307633
$ws.Range("M19").Value = 0.02544333076329992
$ws.Range("N19").Value = 0.0007710100231303007
$ws.Range("O19").Value = 0.06013878180416345
$ws.Range("S19").Value = 0.1187355435620663
